# Append the 2025-05-29 Kaspa buy as new row 29 (A1:D28 -> A1:D29).
#
# Column A holds the date as literal text (e.g. "05/24/2025" in row 28
# is plain text, not a real date value). Assigning a date-shaped string
# straight to a "General" formatted cell's .Value makes Excel parse it
# as an actual date serial, which also forces a brand-new cell style to
# be allocated for that cell (since the original file has zero styled
# cells in that column). Neither happens to match the source workbook.
#
# To get literal text in without Excel's smart date parsing or picking
# up any extra styling, build the text via a formula (TEXT() always
# yields a string, never a date), copy that computed value, and paste
# only the value into the target cell. Paste-values transplants the
# string verbatim instead of re-parsing user input, so no date
# conversion and no new NumberFormat/style gets introduced. The helper
# cell is cleared afterwards so no trace of it remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Cells.Item(31, 1)
$helper.Formula = "=TEXT(DATE(2025,5,29),""mm/dd/yyyy"")"
$ws.Range("A31").Copy()
$ws.Range("A29").PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()

$ws.Cells.Item(29, 2).Value = 508.9639999999999
$ws.Cells.Item(29, 3).Value = 0.09823877523754138
$ws.Cells.Item(29, 4).Value = 50
